# adding new progress as of date 04 nov 2025
#
# For each training row on the "Training Dashboard" sheet (rows 3-30):
#   - decrement column H ("PERIOD TO EXPIRE") by 1 day
#   - set column I ("LAST UPDATE") to the new progress date 04-Nov-2025
#
# Column I holds its date as plain text (not a real date value), so a
# direct Value assignment of a date-looking string would get
# auto-converted into a date serial number and pick up a date number
# format. To avoid that, we write the text via a quoted formula and
# then flatten it back down to a literal value with copy/paste-special,
# which keeps the original cell style and stores a plain text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 30; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H - PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # column I - LAST UPDATE

    $oldPeriod = $hCell.Value2
    $hCell.Value = $oldPeriod - 1

    $iCell.Formula = "=""04-Nov-2025"""
    $iCell.Copy()
    $iCell.PasteSpecial(-4163)         # xlPasteValues
}

$excel.CutCopyMode = 0
